# The "想去人数" (want-to-go count) numbers increased for several events.
# The same events are listed both on the "展览" sheet and again on the
# aggregated "全部类型" sheet (at different row numbers), so both places
# need to be updated.

$wb = $excel.ActiveWorkbook

# Sheet "展览": cell -> new value
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 820
$ws1.Range("F12").Value = 556
$ws1.Range("F14").Value = 13182
$ws1.Range("F18").Value = 5415
$ws1.Range("F19").Value = 5559
$ws1.Range("F20").Value = 22

# Sheet "全部类型": same events, different row numbers
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 820
$ws4.Range("F34").Value = 556
$ws4.Range("F36").Value = 13182
$ws4.Range("F41").Value = 5415
$ws4.Range("F42").Value = 5559
$ws4.Range("F43").Value = 22

$wb.Save()
